# logBook.xlsx - "updated 20july morning entries"
# Adds a new log entry (row 67) for 20-Jul-2022 morning work, copying the
# formatting of the previous entry (row 66) and then filling in the new
# entry's own data. Also updates the current selection to D67, matching
# what Excel records after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles, row height, etc.) of the last existing
# entry row (66) down into the new row (67) before writing the new values.
$ws.Range("A66:G66").Copy($ws.Range("A67:G67"))
$ws.Rows.Item(67).RowHeight = $ws.Rows.Item(66).RowHeight

# New entry: Sno, Date, startTime, endTime, Time(formula), Category, Description
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 44762
$ws.Range("C67").Value = 0.34722222222222227
$ws.Range("D67").Value = 0.375
$ws.Range("E67").Formula = "=D67-C67"
$ws.Range("F67").Value = "Code"
$ws.Range("G67").Value = "1. test video prediction for reference segformer mit b3 model`n2. colab nb 3ep run"

# Recalculate the running total.
$ws.Range("E69").Formula = "=SUM(E2:E68)"

# Update the recorded selection to match the author's final cursor position.
$ws.Range("D67").Select()
